# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values on the zh-cn and de-de sheets to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-30 09:54:20"
$zhcn.Range("H4").Value = "2016-03-30 09:55:06"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-30 09:54:32"
$dede.Range("H4").Value = "2016-03-30 09:55:21"
